$d = $word.ActiveDocument

# 1. Highlight the "Add Cakes to the Database" and "Show Cake Details"
#    Heading2 paragraphs (text run + paragraph mark) in yellow.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($t -eq "Add Cakes to the Database`r") -or ($t -eq "Show Cake Details`r")) {
        $p.Range.Font.HighlightColorIndex = 7  # wdYellow
    }
}

# 2. Change the cached "Page 6 of 6" footer field result to "Page 4 of 6"
$d.Content.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, "4", 2)
